$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("46").Insert()

$ws.Range("A46").Value = 4
$ws.Range("B46").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C46").Value = "Los Lagos"
$ws.Range("D46").Value = 44477
$ws.Range("E46").Value = 10
$ws.Range("F46").Value = 100112039
$ws.Range("G46").Value = "Ciboulette"
$ws.Range("H46").Value = "Sin especificar"
$ws.Range("I46").Value = "Primera"
$ws.Range("J46").Value = 240
$ws.Range("K46").Value = 3000
$ws.Range("L46").Value = 3000
$ws.Range("M46").Value = 3000
$ws.Range("N46").Value = "$/docena de atados"
$ws.Range("O46").Value = "Región Metropolitana"
$ws.Range("P46").Value = 1000
$ws.Range("Q46").Value = 3
$ws.Range("R46").Value = "Hortaliza"
